$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status text (shared string used by E2/F2) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# widen columns E:F on Overview
$wsOverview.Range("E1").ColumnWidth = 29.9777050018311
$wsOverview.Range("F1").ColumnWidth = 29.9777050018311

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# widen columns C, I, J
$wsZh.Range("C1").ColumnWidth = 29.9777050018311
$wsZh.Range("I1").ColumnWidth = 40
$wsZh.Range("J1").ColumnWidth = 40

# Latest Target File (I2) -> hyperlink to md file
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/86f213681a113e07bbbb7fa28e26e1c333635c43/e2e/58ba7908-ddd0-4151-aa88-5daafbd0360coooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$mdDisplay = "58ba7908-ddd0-4151-aa88-5daafbd0360coooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, "", "", $mdDisplay) | Out-Null

# Latest Handback File (J2)
$wsZh.Range("J2").Value = "58ba7908-ddd0-4151-aa88-5daafbd0360cooooooooooooooooooooooooooooooooooooooooo.6339d507da0d5f1fdbeb10305ca8b3d8c760a30a.zh-cn.xlf"

# Latest Handback DateTime (K2)
$wsZh.Range("K2").Value = "2016-10-10 10:09:31"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C1").ColumnWidth = 29.9777050018311
$wsDe.Range("I1").ColumnWidth = 40
$wsDe.Range("J1").ColumnWidth = 40

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, "", "", $mdDisplay) | Out-Null

$wsDe.Range("J2").Value = "58ba7908-ddd0-4151-aa88-5daafbd0360cooooooooooooooooooooooooooooooooooooooooo.6339d507da0d5f1fdbeb10305ca8b3d8c760a30a.de-de.xlf"

$wsDe.Range("K2").Value = "2016-10-10 10:09:47"
